# Delete duplicates function added
# Strips the trailing document-type text (e.g. "답변서", "소장") that
# follows the "(yy.mm.dd)" prefix in column C, for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^(\([^)]*\))') {
        $cell.Value = $matches[1]
    }
}
